# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the per-language report sheets (zh-cn, de-de) for the first
# data row, as would happen when the handback report is regenerated.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 02:58:56"
$wsZhCn.Range("H2").Value = "2016-03-24 02:59:43"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 02:59:04"
$wsDeDe.Range("H2").Value = "2016-03-24 02:59:58"
